# Apply updated cryptocurrency price/volume figures to the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.081.77'
$ws.Range('D3').Value = '1.656.61'
$ws.Range('E3').Value = '  +3.91%  '
$ws.Range('E4').Value = '  +0.07%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '215.25'
$c.ClearFormats()
$ws.Range('E5').Value = '  +1.92%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  +2.00%  '
$ws.Range('E9').Value = '  +1.73%  '
$ws.Range('E10').Value = '  +3.95%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.0866'
$c.ClearFormats()
$ws.Range('E11').Value = '  +1.35%  '
$ws.Range('D12').Value = '1.891.76'
$ws.Range('E12').Value = '  +3.99%  '
$ws.Range('D13').Value = '1.650.49'
$ws.Range('E13').Value = '  +3.42%  '
$ws.Range('E14').Value = '  +2.22%  '
$ws.Range('E15').Value = '  +3.11%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '65.12'
$c.ClearFormats()
$ws.Range('E16').Value = '  +2.57%  '
$ws.Range('D17').Value = '27.072.97'
$ws.Range('E17').Value = '  +3.11%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '238.52'
$c.ClearFormats()
$ws.Range('E18').Value = '  +3.83%  '
$ws.Range('E19').Value = '  +3.32%  '
$ws.Range('E20').Value = '  +1.27%  '
$ws.Range('E21').Value = '  +0.11%  '
$ws.Range('E22').Value = '  +4.56%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '2.25'
$c.ClearFormats()
$ws.Range('E23').Value = '  +4.59%  '
$ws.Range('E24').Value = '  +3.75%  '
$ws.Range('E25').Value = '  -0.22%  '
$ws.Range('E26').Value = '  -0.02%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '7.12'
$c.ClearFormats()
$ws.Range('E27').Value = '  +2.15%  '
$ws.Range('E28').Value = '  +1.29%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '15.83'
$c.ClearFormats()
$ws.Range('E29').Value = '  +3.33%  '
$ws.Range('E30').Value = '  +0.81%  '
$ws.Range('E31').Value = '  +1.89%  '
$ws.Range('E32').Value = '  +3.31%  '
$ws.Range('D33').Value = '1.516.40'
$ws.Range('E33').Value = '  +3.40%  '
$ws.Range('E34').Value = '  +4.46%  '
$ws.Range('E35').Value = '  +10.24%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '2.42'
$c.ClearFormats()
$ws.Range('E36').Value = '  -0.01%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.575'
$c.ClearFormats()
$ws.Range('E37').Value = '  +1.97%  '
$ws.Range('E38').Value = '  +8.58%  '
$ws.Range('E39').Value = '  +2.88%  '
$ws.Range('E40').Value = '  +3.90%  '
$ws.Range('E41').Value = '  +0.05%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '2.26'
$c.ClearFormats()
$ws.Range('E42').Value = '  +4.26%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '66.04'
$c.ClearFormats()
$ws.Range('E43').Value = '  +9.41%  '
$ws.Range('D44').Value = '1.798.01'
$ws.Range('E44').Value = '  +3.70%  '
$ws.Range('E45').Value = '  +3.74%  '
$ws.Range('E46').Value = '  -1.16%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '89.50'
$c.ClearFormats()
$ws.Range('E47').Value = '  +1.97%  '
$ws.Range('E48').Value = '  +2.22%  '
$ws.Range('E49').Value = '  +3.45%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '0.0978'
$c.ClearFormats()
$ws.Range('E51').Value = '  +3.16%  '
